$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (14) of data to the sheet
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 1.0255092592592592
$ws.Range("B14").NumberFormat = "[h]:mm:ss"
$ws.Range("C14").Value = "Matrix (Audiovisual, English, Familiar):38; ¿Quién mató a Sara?  (Subtitled, Spanish, New):39; W.I.T.C.H.  (Audiovisual, English, Familiar):35; Matrix 3 (Audiovisual, English, Familiar):35; Shrek (Audiovisual, English, Familiar):33;"
